# Burndown Sheet update: new Sprint 7 task ("Fix Game Logic" / "Enable
# Flagging" / "Enabe Clicking"), updated estimates/actuals on existing
# tasks, and the Ideal/Actual burndown summary rows + chart shifted down
# by one row to make room for the new task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row above the old "Ideal" summary row (row 10) -----
# This pushes the Ideal/Actual rows (and everything below, including the
# blank filler rows and the chart anchor) down by one.
$ws.Rows.Item(10).Insert()

# The freshly inserted row inherits the merged/"summary" formatting of
# the row that used to be there. Re-stamp it with the plain task-row
# formatting used by rows 4-9 instead (copy format only).
$ws.Range("B9:Q9").Copy()
$ws.Range("B10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Existing task tweaks ---------------------------------------------
$ws.Range("F4").Value = 4
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1

# --- New task: "Fix Game Logic" / "Enable Flagging" (row 9) ----------
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Fix Game Logic"
$ws.Range("D9").Value = "Enable Flagging"
$ws.Range("E9").Value = "Michael"
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2

# --- Continuation of that task: "Enabe Clicking" (new row 10) --------
$ws.Range("D10").Value = "Enabe Clicking"
$ws.Range("E10").Value = "Michael"
$ws.Range("F10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2

# --- Ideal / Actual summary rows now live one row lower (11 / 12) and
# need their top SUM to include the new row 10 of task data -----------
$ws.Range("F11").Formula = "=SUM(F4:F10)"
$ws.Range("F12").Formula = "=SUM(F4:F10)"

# --- Chart series need to follow the Ideal/Actual rows down one ------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$idealSeries = $chart.SeriesCollection(1)
$idealSeries.Formula = "=SERIES(""Ideal Burnout"",,'Burn Down Chart'!`$F`$11:`$M`$11,1)"
$actualSeries = $chart.SeriesCollection(2)
$actualSeries.Formula = "=SERIES(""Actual Burnout"",,'Burn Down Chart'!`$F`$12:`$M`$12,2)"

# The chart floats at a fixed pixel position, so inserting the row above
# doesn't carry it down automatically the way it does for cell-anchored
# content - nudge it down by one (now shorter, 15pt-tall) row so its
# anchor cells land on row 13/row 32 again instead of row 12/row 31.
$co.Top = $co.Top + $ws.Rows.Item(14).Top - $ws.Rows.Item(13).Top

# --- Sprint title / release bump --------------------------------------
$ws.Range("A1").Value = "Project Title: Minesweeper" + [char]10 + "Release #:1.7" + [char]10 + "Sprint #: 7"

# --- Selection cosmetics ----------------------------------------------
$ws.Range("J1").Select()
